$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1).Text = "September 13, 2020"
